$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs, Il13, Il13ra1, ECs)
$ws.Range("I2").Value = 0.9352601111131627
$ws.Range("J2").Value = 0.9352601111131628
$ws.Range("M2").Value = 8.415202000000001
$ws.Range("N2").Value = 25.245606
$ws.Range("O2").Value = 0.1569653516800918
$ws.Range("P2").Value = 0.1569653516800918
$ws.Range("Q2").Value = 2.699655708014
$ws.Range("R2").Value = 24.296901372126
$ws.Range("S2").Value = 0.1468034322532393
$ws.Range("T2").Value = 0.1468034322532393

# Row 3 (FAPs, Il13, Il13ra1, FAPs)
$ws.Range("I3").Value = 0.9352601111131627
$ws.Range("J3").Value = 0.9352601111131628
$ws.Range("O3").Value = 0.5328513631375226
$ws.Range("P3").Value = 0.5328513631375226
$ws.Range("S3").Value = 0.4983546250947996
$ws.Range("T3").Value = 0.4983546250947996

# Row 4 (FAPs, Il13, Il13ra1, MuSCs)
$ws.Range("I4").Value = 0.9352601111131627
$ws.Range("J4").Value = 0.9352601111131628
$ws.Range("M4").Value = 2.790681000000001
$ws.Range("N4").Value = 8.372043000000001
$ws.Range("O4").Value = 0.05205344144940909
$ws.Range("P4").Value = 0.05205344144940908
$ws.Range("Q4").Value = 0.8952699995670003
$ws.Range("R4").Value = 8.057429996103002
$ws.Range("S4").Value = 0.04868350743379685
$ws.Range("T4").Value = 0.04868350743379685

# Row 5 (FAPs, Il13, Il13ra1, Resolving-Mac)
$ws.Range("I5").Value = 0.9352601111131627
$ws.Range("J5").Value = 0.9352601111131628
$ws.Range("M5").Value = 13.838817
$ws.Range("N5").Value = 41.516451
$ws.Range("O5").Value = 0.2581298437329766
$ws.Range("P5").Value = 0.2581298437329766
$ws.Range("Q5").Value = 4.439589365319001
$ws.Range("R5").Value = 39.956304287871
$ws.Range("S5").Value = 0.241418546331327
$ws.Range("T5").Value = 0.241418546331327

# Row 6 (MuSCs, Il13, Il13ra1, ECs)
$ws.Range("G6").Value = 0.02220666666666667
$ws.Range("H6").Value = 0.06662
$ws.Range("I6").Value = 0.06473988888683736
$ws.Range("J6").Value = 0.06473988888683736
$ws.Range("M6").Value = 8.415202000000001
$ws.Range("N6").Value = 25.245606
$ws.Range("O6").Value = 0.1569653516800918
$ws.Range("P6").Value = 0.1569653516800918
$ws.Range("Q6").Value = 0.1868735857466667
$ws.Range("R6").Value = 1.68186227172
$ws.Range("S6").Value = 0.01016191942685249
$ws.Range("T6").Value = 0.01016191942685249

# Row 7 (MuSCs, Il13, Il13ra1, FAPs)
$ws.Range("G7").Value = 0.02220666666666667
$ws.Range("H7").Value = 0.06662
$ws.Range("I7").Value = 0.06473988888683736
$ws.Range("J7").Value = 0.06473988888683736
$ws.Range("O7").Value = 0.5328513631375226
$ws.Range("P7").Value = 0.5328513631375226
$ws.Range("Q7").Value = 0.6343810518288888
$ws.Range("R7").Value = 5.70942946646
$ws.Range("S7").Value = 0.03449673804272304
$ws.Range("T7").Value = 0.03449673804272304

# Row 8 (MuSCs, Il13, Il13ra1, MuSCs)
$ws.Range("G8").Value = 0.02220666666666667
$ws.Range("H8").Value = 0.06662
$ws.Range("I8").Value = 0.06473988888683736
$ws.Range("J8").Value = 0.06473988888683736
$ws.Range("M8").Value = 2.790681000000001
$ws.Range("N8").Value = 8.372043000000001
$ws.Range("O8").Value = 0.05205344144940909
$ws.Range("P8").Value = 0.05205344144940908
$ws.Range("Q8").Value = 0.06197172274000001
$ws.Range("R8").Value = 0.5577455046600001
$ws.Range("S8").Value = 0.003369934015612239
$ws.Range("T8").Value = 0.003369934015612238

# Row 9 (MuSCs, Il13, Il13ra1, Resolving-Mac)
$ws.Range("G9").Value = 0.02220666666666667
$ws.Range("H9").Value = 0.06662
$ws.Range("I9").Value = 0.06473988888683736
$ws.Range("J9").Value = 0.06473988888683736
$ws.Range("M9").Value = 13.838817
$ws.Range("N9").Value = 41.516451
$ws.Range("O9").Value = 0.2581298437329766
$ws.Range("P9").Value = 0.2581298437329766
$ws.Range("Q9").Value = 0.30731399618
$ws.Range("R9").Value = 2.76582596562
$ws.Range("S9").Value = 0.0167112974016496
$ws.Range("T9").Value = 0.0167112974016496
